# Add the new docgen variable "docgen.bookmark.always_bold" as row 16 on
# both the "en" and "ko" worksheets of the docgen_variables workbook.
#
# Column order on every existing data row is:
#   A = variable name   B = type ("boolean"/"string"/"int"/"float"/"table")
#   C = localized description (wrapped text)
#
# NOTE on shared-string ordering: the canonical file stores the new
# Korean description ahead of the new English description in
# xl/sharedStrings.xml (the variable-name string is common to both
# sheets and is appended first). Writing the "ko" sheet's C16 value
# before the "en" sheet's C16 value reproduces that exact ordering.

$wb = $excel.ActiveWorkbook

# --- Korean sheet ("ko") ---
$wsKo = $wb.Worksheets.Item("ko")
$wsKo.Range("A16").Value = "docgen.bookmark.always_bold"
$wsKo.Range("B16").Value = "boolean"
$wsKo.Range("B16").HorizontalAlignment = -4108   # xlCenter
$wsKo.Range("C16").Value = "blookmark 참조를 항상 굵은 글씨로 표현합니다.`n(기본값 : false)"
$wsKo.Range("C16").WrapText = $true
$wsKo.Rows.Item(16).RowHeight = 33

# --- English sheet ("en") ---
$wsEn = $wb.Worksheets.Item("en")
$wsEn.Range("A16").Value = "docgen.bookmark.always_bold"
$wsEn.Range("B16").Value = "boolean"
$wsEn.Range("B16").HorizontalAlignment = -4108   # xlCenter
$wsEn.Range("C16").Value = "blookmark text always display in bold.`n(default : false)"
$wsEn.Range("C16").WrapText = $true
$wsEn.Rows.Item(16).RowHeight = 33

# --- Restore selection state: both sheets end up with A16 selected, and
#     "ko" remains the active/tabbed sheet (as in the source workbook). ---
$wsEn.Activate()
$wsEn.Range("A16").Select()
$wsKo.Activate()
$wsKo.Range("A16").Select()
